$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 = 1, formatted the same as B1 (bold font, thin border, centered/top)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = 1

# New column C data rows 2-25, all 200 (plain, unstyled like column B data cells)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 200
}
